$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.764.42"
$ws.Range("E2").Value = "  -2.36%  "
$ws.Range("D3").Value = "1.563.93"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.19"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("E6").Value = "  -2.71%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0861"
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("D12").Value = "1.786.11"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "1.564.54"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.72"
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.513"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").Value = "26.813.26"
$ws.Range("E16").Value = "  -2.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.51"
$ws.Range("E17").Value = "  -3.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "213.90"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.34"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.08"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.31"
$ws.Range("E23").Value = "  -2.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.09"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.90"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0463"
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.11"
$ws.Range("E31").Value = "  -3.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.14"
$ws.Range("E32").Value = "  -1.96%  "
$ws.Range("D33").Value = "1.382.12"
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.92"
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.924"
$ws.Range("E37").Value = "  -4.31%  "
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.522"
$ws.Range("E39").Value = "  -1.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.813"
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.993"
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.36"
$ws.Range("E43").Value = "  +1.61%  "
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.20"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").Value = "1.699.24"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.41"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("D49").Value = "0.0₇0986"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0949"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("E51").Value = "  -0.69%  "
